$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.0606410827001
$ws.Cells.Item(2, 4).Value = 1.070574647351912
$ws.Cells.Item(2, 5).Value = 1.066926387837632
$ws.Cells.Item(2, 6).Value = 1.078938631139831
$ws.Cells.Item(2, 9).Value = 1.039818476757014
$ws.Cells.Item(2, 10).Value = 1.06562088253957
$ws.Cells.Item(2, 11).Value = 1.073274120917788
$ws.Cells.Item(2, 12).Value = 1.069635624528051
$ws.Cells.Item(2, 13).Value = 1.08161599103682
$ws.Cells.Item(2, 14).Value = 1.025497224839458

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.061957523776998
$ws.Cells.Item(3, 4).Value = 1.071834941636991
$ws.Cells.Item(3, 5).Value = 1.068118520191134
$ws.Cells.Item(3, 6).Value = 1.08026371582058
$ws.Cells.Item(3, 9).Value = 1.040030688710887
$ws.Cells.Item(3, 10).Value = 1.066589301096673
$ws.Cells.Item(3, 11).Value = 1.074349705372908
$ws.Cells.Item(3, 12).Value = 1.070642505546355
$ws.Cells.Item(3, 13).Value = 1.082757817961417
$ws.Cells.Item(3, 14).Value = 1.025828041571685

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.062808820180296
$ws.Cells.Item(4, 4).Value = 1.072650234621578
$ws.Cells.Item(4, 5).Value = 1.068889752340983
$ws.Cells.Item(4, 6).Value = 1.081121141640002
$ws.Cells.Item(4, 9).Value = 1.040166259678381
$ws.Cells.Item(4, 10).Value = 1.067214900441534
$ws.Cells.Item(4, 11).Value = 1.075044916158727
$ws.Cells.Item(4, 12).Value = 1.071293292057346
$ws.Cells.Item(4, 13).Value = 1.08349611116963
$ws.Cells.Item(4, 14).Value = 1.026041548570605

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.063166581838047
$ws.Cells.Item(5, 4).Value = 1.072992938063536
$ws.Cells.Item(5, 5).Value = 1.069213943002508
$ws.Cells.Item(5, 6).Value = 1.081481607924022
$ws.Cells.Item(5, 9).Value = 1.040222836582019
$ws.Cells.Item(5, 10).Value = 1.067477657424021
$ws.Cells.Item(5, 11).Value = 1.075337002106305
$ws.Cells.Item(5, 12).Value = 1.071566709384516
$ws.Cells.Item(5, 13).Value = 1.083806361522915
$ws.Cells.Item(5, 14).Value = 1.026131175008608

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.063226644456891
$ws.Cells.Item(6, 4).Value = 1.073050476872015
$ws.Cells.Item(6, 5).Value = 1.069268374035759
$ws.Cells.Item(6, 6).Value = 1.081542132139438
$ws.Cells.Item(6, 9).Value = 1.04023231164596
$ws.Cells.Item(6, 10).Value = 1.067521761173553
$ws.Cells.Item(6, 11).Value = 1.075386034094555
$ws.Cells.Item(6, 12).Value = 1.071612607254647
$ws.Cells.Item(6, 13).Value = 1.083858446479552
$ws.Cells.Item(6, 14).Value = 1.026146215964834

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.062813601089661
$ws.Cells.Item(7, 4).Value = 1.072654814018249
$ws.Cells.Item(7, 5).Value = 1.068894084328454
$ws.Cells.Item(7, 6).Value = 1.081125958187223
$ws.Cells.Item(7, 9).Value = 1.040167017300196
$ws.Cells.Item(7, 10).Value = 1.067218412373164
$ws.Cells.Item(7, 11).Value = 1.075048819732495
$ws.Cells.Item(7, 12).Value = 1.071296946151462
$ws.Cells.Item(7, 13).Value = 1.083500257250446
$ws.Cells.Item(7, 14).Value = 1.026042746680987

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.06108609164373
$ws.Cells.Item(8, 4).Value = 1.07100061300909
$ws.Cells.Item(8, 5).Value = 1.06732930848352
$ws.Cells.Item(8, 6).Value = 1.079386449797063
$ws.Cells.Item(8, 9).Value = 1.039890556246009
$ws.Cells.Item(8, 10).Value = 1.065948379138838
$ws.Cells.Item(8, 11).Value = 1.073637779357987
$ws.Cells.Item(8, 12).Value = 1.069976057699553
$ws.Cells.Item(8, 13).Value = 1.082001991119899
$ws.Cells.Item(8, 14).Value = 1.025609140890756

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.0580377999081
$ws.Cells.Item(9, 4).Value = 1.068084039768255
$ws.Cells.Item(9, 5).Value = 1.064570666071363
$ws.Cells.Item(9, 6).Value = 1.07632114959655
$ws.Cells.Item(9, 9).Value = 1.039390015398873
$ws.Cells.Item(9, 10).Value = 1.063702416742467
$ws.Cells.Item(9, 11).Value = 1.071145388832206
$ws.Cells.Item(9, 12).Value = 1.067642775391619
$ws.Cells.Item(9, 13).Value = 1.079357571413789
$ws.Cells.Item(9, 14).Value = 1.024840804710559

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.056002573140467
$ws.Cells.Item(10, 4).Value = 1.066138378430726
$ws.Cells.Item(10, 5).Value = 1.06273054489967
$ws.Cells.Item(10, 6).Value = 1.074277401279447
$ws.Cells.Item(10, 9).Value = 1.039047293499855
$ws.Cells.Item(10, 10).Value = 1.062199592193636
$ws.Cells.Item(10, 11).Value = 1.069479642691176
$ws.Cells.Item(10, 12).Value = 1.066083280377605
$ws.Cells.Item(10, 13).Value = 1.077591591613567
$ws.Cells.Item(10, 14).Value = 1.024325671466603

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.055120531129342
$ws.Cells.Item(11, 4).Value = 1.065295545081883
$ws.Cells.Item(11, 5).Value = 1.061933475194745
$ws.Cells.Item(11, 6).Value = 1.073392345759836
$ws.Cells.Item(11, 9).Value = 1.038896742442051
$ws.Cells.Item(11, 10).Value = 1.061547513114426
$ws.Cells.Item(11, 11).Value = 1.068757337824279
$ws.Cells.Item(11, 12).Value = 1.065407028477327
$ws.Cells.Item(11, 13).Value = 1.076826148518496
$ws.Cells.Item(11, 14).Value = 1.024101913862449

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.054792780806825
$ws.Cells.Item(12, 4).Value = 1.064982423499641
$ws.Cells.Item(12, 5).Value = 1.061637362203196
$ws.Cells.Item(12, 6).Value = 1.073063578404667
$ws.Cells.Item(12, 9).Value = 1.038840497331169
$ws.Cells.Item(12, 10).Value = 1.06130509689753
$ws.Cells.Item(12, 11).Value = 1.068488884540296
$ws.Cells.Item(12, 12).Value = 1.065155688384521
$ws.Cells.Item(12, 13).Value = 1.076541711260957
$ws.Cells.Item(12, 14).Value = 1.024018694095158

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.054863089889968
$ws.Cells.Item(13, 4).Value = 1.065049591717781
$ws.Cells.Item(13, 5).Value = 1.061700881566263
$ws.Cells.Item(13, 6).Value = 1.073134101003381
$ws.Cells.Item(13, 9).Value = 1.038852576759797
$ws.Cells.Item(13, 10).Value = 1.061357105320886
$ws.Cells.Item(13, 11).Value = 1.068546475837947
$ws.Cells.Item(13, 12).Value = 1.065209608521007
$ws.Cells.Item(13, 13).Value = 1.076602729390325
$ws.Cells.Item(13, 14).Value = 1.024036549838373

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.055093441670679
$ws.Cells.Item(14, 4).Value = 1.065269663510593
$ws.Cells.Item(14, 5).Value = 1.061908999354373
$ws.Cells.Item(14, 6).Value = 1.073365170141614
$ws.Cells.Item(14, 9).Value = 1.038892099811997
$ws.Cells.Item(14, 10).Value = 1.061527479112602
$ws.Cells.Item(14, 11).Value = 1.068735150614086
$ws.Cells.Item(14, 12).Value = 1.065386255703325
$ws.Cells.Item(14, 13).Value = 1.076802639263267
$ws.Cells.Item(14, 14).Value = 1.024095037062447

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.05523535289396
$ws.Cells.Item(15, 4).Value = 1.065405249519308
$ws.Cells.Item(15, 5).Value = 1.062037221432189
$ws.Cells.Item(15, 6).Value = 1.073507536931711
$ws.Cells.Item(15, 9).Value = 1.038916408348558
$ws.Cells.Item(15, 10).Value = 1.061632424779348
$ws.Cells.Item(15, 11).Value = 1.068851378464558
$ws.Cells.Item(15, 12).Value = 1.06549507389984
$ws.Cells.Item(15, 13).Value = 1.076925794656283
$ws.Cells.Item(15, 14).Value = 1.024131058869863

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.056061094673141
$ws.Cells.Item(16, 4).Value = 1.066194306867345
$ws.Cells.Item(16, 5).Value = 1.062783437530159
$ws.Cells.Item(16, 6).Value = 1.074336137119031
$ws.Cells.Item(16, 9).Value = 1.039057239715616
$ws.Cells.Item(16, 10).Value = 1.062242839936814
$ws.Cells.Item(16, 11).Value = 1.069527557819846
$ws.Cells.Item(16, 12).Value = 1.066128140101279
$ws.Cells.Item(16, 13).Value = 1.077642375228134
$ws.Cells.Item(16, 14).Value = 1.024340506670593

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.056578850154322
$ws.Cells.Item(17, 4).Value = 1.066689165897109
$ws.Cells.Item(17, 5).Value = 1.063251440745199
$ws.Cells.Item(17, 6).Value = 1.074855867184014
$ws.Cells.Item(17, 9).Value = 1.039145003345775
$ws.Cells.Item(17, 10).Value = 1.062625374887286
$ws.Cells.Item(17, 11).Value = 1.069951430529569
$ws.Cells.Item(17, 12).Value = 1.066524981418957
$ws.Cells.Item(17, 13).Value = 1.078091661044999
$ws.Cells.Item(17, 14).Value = 1.024471699339413

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.056880773635306
$ws.Cells.Item(18, 4).Value = 1.066977775333301
$ws.Cells.Item(18, 5).Value = 1.063524391798314
$ws.Cells.Item(18, 6).Value = 1.07515900761981
$ws.Cells.Item(18, 9).Value = 1.039195986969884
$ws.Cells.Item(18, 10).Value = 1.062848371358565
$ws.Cells.Item(18, 11).Value = 1.070198569538573
$ws.Cells.Item(18, 12).Value = 1.066756357851499
$ws.Cells.Item(18, 13).Value = 1.078353648531849
$ws.Cells.Item(18, 14).Value = 1.024548154175148

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.056983709223032
$ws.Cells.Item(19, 4).Value = 1.067076178115645
$ws.Cells.Item(19, 5).Value = 1.063617456511988
$ws.Cells.Item(19, 6).Value = 1.075262369197813
$ws.Cells.Item(19, 9).Value = 1.03921333590764
$ws.Cells.Item(19, 10).Value = 1.062924385477539
$ws.Cells.Item(19, 11).Value = 1.070282820890249
$ws.Cells.Item(19, 12).Value = 1.0668352353075
$ws.Cells.Item(19, 13).Value = 1.078442967149806
$ws.Cells.Item(19, 14).Value = 1.024574211852478

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.05652330763178
$ws.Cells.Item(20, 4).Value = 1.066636075692907
$ws.Cells.Item(20, 5).Value = 1.063201231276332
$ws.Cells.Item(20, 6).Value = 1.074800106056821
$ws.Cells.Item(20, 9).Value = 1.039135608598342
$ws.Cells.Item(20, 10).Value = 1.062584345943302
$ws.Cells.Item(20, 11).Value = 1.069905963238408
$ws.Cells.Item(20, 12).Value = 1.066482413885855
$ws.Cells.Item(20, 13).Value = 1.078043464537223
$ws.Cells.Item(20, 14).Value = 1.024457630609189

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.055025612176505
$ws.Cells.Item(21, 4).Value = 1.06520485941824
$ws.Cells.Item(21, 5).Value = 1.061847715141674
$ws.Cells.Item(21, 6).Value = 1.073297126573407
$ws.Cells.Item(21, 9).Value = 1.038880470200792
$ws.Cells.Item(21, 10).Value = 1.061477313959091
$ws.Cells.Item(21, 11).Value = 1.068679594955739
$ws.Cells.Item(21, 12).Value = 1.065334241674601
$ws.Cells.Item(21, 13).Value = 1.076743774059914
$ws.Cells.Item(21, 14).Value = 1.024077816971934

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.054083251507021
$ws.Cells.Item(22, 4).Value = 1.064304672750608
$ws.Cells.Item(22, 5).Value = 1.06099643927612
$ws.Cells.Item(22, 6).Value = 1.072352035956856
$ws.Cells.Item(22, 9).Value = 1.03871818112672
$ws.Cells.Item(22, 10).Value = 1.060780091487683
$ws.Cells.Item(22, 11).Value = 1.067907618147116
$ws.Cells.Item(22, 12).Value = 1.064611470851922
$ws.Cells.Item(22, 13).Value = 1.075925925066907
$ws.Cells.Item(22, 14).Value = 1.023838398029502

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.054582882479006
$ws.Cells.Item(23, 4).Value = 1.064781910567793
$ws.Cells.Item(23, 5).Value = 1.061447742983705
$ws.Cells.Item(23, 6).Value = 1.072853057745253
$ws.Cells.Item(23, 9).Value = 1.038804391450264
$ws.Cells.Item(23, 10).Value = 1.061149815803738
$ws.Cells.Item(23, 11).Value = 1.068316944953185
$ws.Cells.Item(23, 12).Value = 1.064994708669742
$ws.Cells.Item(23, 13).Value = 1.076359547976947
$ws.Cells.Item(23, 14).Value = 1.023965377119397

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.056548405137042
$ws.Cells.Item(24, 4).Value = 1.066660064974262
$ws.Cells.Item(24, 5).Value = 1.063223918858378
$ws.Cells.Item(24, 6).Value = 1.074825302138695
$ws.Cells.Item(24, 9).Value = 1.039139854321201
$ws.Cells.Item(24, 10).Value = 1.062602885557672
$ws.Cells.Item(24, 11).Value = 1.069926508255274
$ws.Cells.Item(24, 12).Value = 1.066501648613358
$ws.Cells.Item(24, 13).Value = 1.078065242690986
$ws.Cells.Item(24, 14).Value = 1.02446398787175

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.058826375489966
$ws.Cells.Item(25, 4).Value = 1.068838257307755
$ws.Cells.Item(25, 5).Value = 1.065284011304891
$ws.Cells.Item(25, 6).Value = 1.077113628407384
$ws.Cells.Item(25, 9).Value = 1.039521006142723
$ws.Cells.Item(25, 10).Value = 1.064284014093119
$ws.Cells.Item(25, 11).Value = 1.071790452586833
$ws.Cells.Item(25, 12).Value = 1.068246675161169
$ws.Cells.Item(25, 13).Value = 1.080041740646256
$ws.Cells.Item(25, 14).Value = 1.025039947834108

